$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 382.26086
$ws.Range("I6").Value = 382.26086
$ws.Range("K6").Value = 1146.78258
$ws.Range("M6").Value = -1034.78258
$ws.Range("H28").Value = 1097.2858
$ws.Range("I28").Value = 929.44446
$ws.Range("K28").Value = 929.44446
$ws.Range("M28").Value = -444.44446
$ws.Range("H32").Value = 11114500
$ws.Range("I32").Value = 1000
$ws.Range("K32").Value = 1000
$ws.Range("M32").Value = -674
$ws.Range("H40").Value = 4882.72
$ws.Range("I40").Value = 5530.4
$ws.Range("J40").Value = 3911.2
$ws.Range("K40").Value = 5530.4
$ws.Range("L40").Value = 3911.2
$ws.Range("M40").Value = -5355.4
$ws.Range("N40").Value = -4261.2
$ws.Range("H41").Value = 999.375
$ws.Range("I41").Value = 856.4286
$ws.Range("K41").Value = 856.4286
$ws.Range("M41").Value = -416.4286
$ws.Range("H62").Value = 13526.909
$ws.Range("I62").Value = 18399.857
$ws.Range("J62").Value = 4999.25
$ws.Range("K62").Value = 18399.857
$ws.Range("L62").Value = 4999.25
$ws.Range("M62").Value = -17775.857
$ws.Range("N62").Value = -6247.25
$ws.Range("H65").Value = 13526.909
$ws.Range("I65").Value = 18399.857
$ws.Range("J65").Value = 4999.25
$ws.Range("K65").Value = 91999.285
$ws.Range("L65").Value = 24996.25
$ws.Range("M65").Value = -88879.285
$ws.Range("N65").Value = -31236.25
$ws.Range("H70").Value = 8134.7
$ws.Range("I70").Value = 10499
$ws.Range("J70").Value = 7543.625
$ws.Range("K70").Value = 31497
$ws.Range("L70").Value = 22630.875
$ws.Range("M70").Value = -31227
$ws.Range("N70").Value = -23170.875
$ws.Range("H73").Value = 8134.7
$ws.Range("I73").Value = 10499
$ws.Range("J73").Value = 7543.625
$ws.Range("K73").Value = 31497
$ws.Range("L73").Value = 22630.875
$ws.Range("M73").Value = -30561
$ws.Range("N73").Value = -24502.875
$ws.Range("H76").Value = 4439.6
$ws.Range("I76").Value = 4424.5
$ws.Range("J76").Value = 4500
$ws.Range("K76").Value = 4424.5
$ws.Range("L76").Value = 4500
$ws.Range("M76").Value = -4109.5
$ws.Range("N76").Value = -5130
$ws.Range("H79").Value = 4439.6
$ws.Range("I79").Value = 4424.5
$ws.Range("J79").Value = 4500
$ws.Range("K79").Value = 4424.5
$ws.Range("L79").Value = 4500
$ws.Range("M79").Value = -3332.5
$ws.Range("N79").Value = -6684
$ws.Range("H130").Value = 20000
$ws.Range("J130").Value = 20000
$ws.Range("L130").Value = 20000
$ws.Range("N130").Value = -30040
$ws.Range("H141").Value = 6473.732
$ws.Range("I141").Value = 4511.75
$ws.Range("K141").Value = 13535.25
$ws.Range("M141").Value = -8355.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3335.8438
$ws.Range("I2").Value = 3025.88
$ws.Range("J2").Value = 4442.857
$ws.Range("K2").Value = 3025.88
$ws.Range("L2").Value = 4442.857
$ws.Range("M2").Value = -2912.88
$ws.Range("N2").Value = -4668.857
$ws.Range("H34").Value = 145006.75
$ws.Range("J34").Value = 145006.75
$ws.Range("L34").Value = 145006.75
$ws.Range("N34").Value = -145548.75
$ws.Range("H61").Value = 4184.8
$ws.Range("I61").Value = 4184.8
$ws.Range("K61").Value = 4184.8
$ws.Range("M61").Value = -3972.8
$ws.Range("H116").Value = 3335.8438
$ws.Range("I116").Value = 3025.88
$ws.Range("J116").Value = 4442.857
$ws.Range("K116").Value = 3025.88
$ws.Range("L116").Value = 4442.857
$ws.Range("M116").Value = -731.8800000000001
$ws.Range("N116").Value = -9030.857
$ws.Range("H136").Value = 4184.8
$ws.Range("I136").Value = 4184.8
$ws.Range("K136").Value = 12554.4
$ws.Range("M136").Value = -10004.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3335.8438
$ws.Range("I3").Value = 3025.88
$ws.Range("J3").Value = 4442.857
$ws.Range("K3").Value = 3025.88
$ws.Range("L3").Value = 4442.857
$ws.Range("M3").Value = -2911.88
$ws.Range("N3").Value = -4670.857
$ws.Range("H51").Value = 79979.5
$ws.Range("J51").Value = 79979.5
$ws.Range("L51").Value = 79979.5
$ws.Range("N51").Value = -80961.5
$ws.Range("H94").Value = 3219.6
$ws.Range("J94").Value = 3452
$ws.Range("L94").Value = 3452
$ws.Range("N94").Value = -4354
$ws.Range("H99").Value = 4304.625
$ws.Range("I99").Value = 4222
$ws.Range("J99").Value = 4332.1665
$ws.Range("K99").Value = 4222
$ws.Range("L99").Value = 4332.1665
$ws.Range("M99").Value = -2724
$ws.Range("N99").Value = -7328.1665
$ws.Range("H103").Value = 26303
$ws.Range("J103").Value = 26303
$ws.Range("L103").Value = 26303
$ws.Range("N103").Value = -28647

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 636.3570999999999
$ws.Range("I19").Value = 678.2727
$ws.Range("J19").Value = 482.66666
$ws.Range("K19").Value = 678.2727
$ws.Range("L19").Value = 482.66666
$ws.Range("M19").Value = -508.2727
$ws.Range("N19").Value = -822.66666
$ws.Range("H24").Value = 636.3570999999999
$ws.Range("I24").Value = 678.2727
$ws.Range("J24").Value = 482.66666
$ws.Range("K24").Value = 678.2727
$ws.Range("L24").Value = 482.66666
$ws.Range("M24").Value = -508.2727
$ws.Range("N24").Value = -822.66666
$ws.Range("H99").Value = 14966284
$ws.Range("I99").Value = 2223023
$ws.Range("K99").Value = 2223023
$ws.Range("M99").Value = -2221525
$ws.Range("H105").Value = 1599.2106
$ws.Range("I105").Value = 1538
$ws.Range("K105").Value = 1538
$ws.Range("M105").Value = 209
$ws.Range("H126").Value = 14966284
$ws.Range("I126").Value = 2223023
$ws.Range("K126").Value = 6669069
$ws.Range("M126").Value = -6666599

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 797.8
$ws.Range("I25").Value = 63
$ws.Range("J25").Value = 1900
$ws.Range("K25").Value = 189
$ws.Range("L25").Value = 5700
$ws.Range("M25").Value = -20
$ws.Range("N25").Value = -6038
$ws.Range("H30").Value = 797.8
$ws.Range("I30").Value = 63
$ws.Range("J30").Value = 1900
$ws.Range("K30").Value = 189
$ws.Range("L30").Value = 5700
$ws.Range("M30").Value = -87
$ws.Range("N30").Value = -5904
$ws.Range("H38").Value = 1060.2
$ws.Range("I38").Value = 71.31579000000001
$ws.Range("J38").Value = 4191.6665
$ws.Range("K38").Value = 213.94737
$ws.Range("L38").Value = 12574.9995
$ws.Range("M38").Value = 133.05263
$ws.Range("N38").Value = -13268.9995
$ws.Range("H125").Value = 2000
$ws.Range("I125").Value = 2000
$ws.Range("K125").Value = 6000
$ws.Range("M125").Value = -1080
$ws.Range("H129").Value = 1528.5
$ws.Range("I129").Value = 1032.125
$ws.Range("J129").Value = 2024.875
$ws.Range("K129").Value = 3096.375
$ws.Range("L129").Value = 6074.625
$ws.Range("M129").Value = 1903.625
$ws.Range("N129").Value = -16074.625
$ws.Range("H137").Value = 1675
$ws.Range("I137").Value = 1675
$ws.Range("K137").Value = 5025
$ws.Range("M137").Value = 75
$ws.Range("H139").Value = 142861440
$ws.Range("I139").Value = 200003410
$ws.Range("J139").Value = 6500
$ws.Range("K139").Value = 600010230
$ws.Range("L139").Value = 19500
$ws.Range("M139").Value = -600005090
$ws.Range("N139").Value = -29780
$ws.Range("H140").Value = 86668330
$ws.Range("I140").Value = 86668330
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 260004990
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -259999810
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1075
$ws.Range("I31").Value = 766.6667
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 766.6667
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = -474.6667
$ws.Range("N31").Value = -2584
$ws.Range("H37").Value = 1075
$ws.Range("I37").Value = 766.6667
$ws.Range("J37").Value = 2000
$ws.Range("K37").Value = 766.6667
$ws.Range("L37").Value = 2000
$ws.Range("M37").Value = -489.6667
$ws.Range("N37").Value = -2554
$ws.Range("H38").Value = 16598.4
$ws.Range("J38").Value = 16598.4
$ws.Range("L38").Value = 16598.4
$ws.Range("N38").Value = -17524.4
$ws.Range("H46").Value = 22513.666
$ws.Range("I46").Value = 22513.666
$ws.Range("K46").Value = 22513.666
$ws.Range("M46").Value = -22357.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 1446.8334
$ws.Range("I11").Value = 1202.5
$ws.Range("K11").Value = 1202.5
$ws.Range("M11").Value = -1062.5
$ws.Range("H122").Value = 4064.4614
$ws.Range("I122").Value = 2483.4
$ws.Range("J122").Value = 9334.666999999999
$ws.Range("K122").Value = 7450.200000000001
$ws.Range("L122").Value = 28004.001
$ws.Range("M122").Value = -5000.200000000001
$ws.Range("N122").Value = -32904.001
$ws.Range("H136").Value = 4968.6665
$ws.Range("I136").Value = 8179.5
$ws.Range("K136").Value = 24538.5
$ws.Range("M136").Value = -21988.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 24999.666
$ws.Range("J40").Value = 29999.5
$ws.Range("L40").Value = 29999.5
$ws.Range("N40").Value = -30297.5
$ws.Range("H96").Value = 40298.355
$ws.Range("I96").Value = 59130.777
$ws.Range("J96").Value = 6400
$ws.Range("K96").Value = 6400
$ws.Range("L96").Value = 6400
$ws.Range("M96").Value = -57757.777
$ws.Range("N96").Value = -9146
